# Apply invoice field updates (first draft of drag drop BO changes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Client / header block ---
$ws.Range("C12").Value = "Yazle Marketing Management"

# Invoice Date (text that looks like a date -> force literal text with a
# leading quote so Excel keeps it as a string instead of coercing it to a
# date serial number, matching the original inlineStr storage).
$ws.Range("F12").Value = "'03/02/2026"

$ws.Range("C13").Value = "Napptix test Address"

# Due Date (same text-forcing trick as above)
$ws.Range("F13").Value = "'05/03/2026"

$ws.Range("C14").Value = "Napptix test Address"

$ws.Range("C15").Value = "Napptix test Address"
$ws.Range("F15").Value = "PD25|22041|4"

# Client VAT No. - force text so the long numeric string is not coerced
# into a number.
$ws.Range("C16").Value = "'100041433200003"

# --- Line item ---
$ws.Range("C21").Value = "Campaign Name"
$ws.Range("D21").Value = 2025
$ws.Range("E21").Value = 14
$ws.Range("F21").Value = 28.35

# --- Subtotal ---
$ws.Range("F25").Value = 28.35

# --- Totals / VAT row ---
$ws.Range("C26").Value = "TWENTY NINE DOLLARS AND SEVENTY SEVEN CENTS"
$ws.Range("E26").Value = "VAT(5%)"
$ws.Range("F26").Value = 1.4175
